$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.438.80"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +5.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.816.67"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +5.08%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9979"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.06"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9977"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5682"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +17.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3850"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +10.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.49"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07645"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +5.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.140"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +8.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.39"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +7.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.9976"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.249"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +6.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.809.35"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +4.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.276"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +6.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.46"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +6.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001083"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +4.85%  "
$ws.Range("E19").Value = "  +1.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9974"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.35"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.009"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +5.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.448.09"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +5.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.34"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.088"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.89"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +4.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.39"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.403"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +16.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.019.73"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +4.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.74"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.158"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +11.45%  "
$ws.Range("E32").Value = "  +12.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.791"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +7.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.636"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02318"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +6.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2160"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +8.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.752"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +16.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "11.71"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +6.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6463"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +8.06%  "
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.066"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +6.55%  "
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06101"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9971"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("E43").Value = "  +3.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.374"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.48"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +5.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6037"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +7.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.710"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.55"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.947"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +5.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.148"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +4.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06853"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.99%  "
